$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "96.641.19"
$ws.Range("E2").Value = "  -1.87%  "

# Row 3
$ws.Range("D3").Value = "3.672.42"
$ws.Range("E3").Value = "  +0.99%  "

# Row 4
$ws.Range("E4").Value = "  -0.11%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.54"
$ws.Range("E5").Value = "  -2.27%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.86"
$ws.Range("E6").Value = "  +6.25%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "654.86"
$ws.Range("E7").Value = "  -0.71%  "

# Row 8
$ws.Range("E8").Value = "  -0.22%  "

# Row 9
$ws.Range("E9").Value = "  -0.55%  "

# Row 10
$ws.Range("E10").Value = "  +0.01%  "

# Row 11
$ws.Range("D11").Value = "3.670.63"
$ws.Range("E11").Value = "  +1.09%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "45.61"
$ws.Range("E12").Value = "  +3.17%  "

# Row 13
$ws.Range("E13").Value = "  -0.69%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.80"
$ws.Range("E14").Value = "  +4.83%  "

# Row 15
$ws.Range("D15").Value = "4.355.03"
$ws.Range("E15").Value = "  +0.91%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000267"
$ws.Range("E16").Value = "  +2.13%  "

# Row 17
$ws.Range("D17").Value = "96.354.26"
$ws.Range("E17").Value = "  -1.71%  "

# Row 18
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.664.26"
$ws.Range("E18").Value = "  +0.77%  "

# Row 19
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "18.70"
$ws.Range("E19").Value = "  +2.42%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.79"
$ws.Range("E20").Value = "  -0.78%  "

# Row 21
$ws.Range("B21").Value = "Polkadot"
$ws.Range("C21").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.78"
$ws.Range("E21").Value = "  -4.94%  "

# Row 22
$ws.Range("E22").Value = "  -2.58%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "524.54"
$ws.Range("E23").Value = "  +1.19%  "

# Row 24
$ws.Range("E24").Value = "  -1.28%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.08"
$ws.Range("E25").Value = "  +1.96%  "

# Row 26
$ws.Range("E26").Value = "  -2.65%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "101.46"
$ws.Range("E27").Value = "  +1.47%  "

# Row 28
$ws.Range("E28").Value = "  +0.87%  "

# Row 29
$ws.Range("D29").Value = "3.867.00"
$ws.Range("E29").Value = "  +0.85%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.168"
$ws.Range("E30").Value = "  +6.71%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.46"
$ws.Range("E31").Value = "  +4.53%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.04"
$ws.Range("E32").Value = "  -1.12%  "

# Row 33
$ws.Range("E33").Value = "  +0.02%  "

# Row 34
$ws.Range("E34").Value = "  +14.94%  "

# Row 35
$ws.Range("E35").Value = "  -1.21%  "

# Row 36
$ws.Range("B36").Value = "Bittensor"
$ws.Range("C36").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "675.13"
$ws.Range("E36").Value = "  +9.54%  "

# Row 37
$ws.Range("B37").Value = "Binance-PegBSC-USD"
$ws.Range("C37").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.61%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "32.39"
$ws.Range("E38").Value = "  +0.87%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.587"
$ws.Range("E39").Value = "  +1.88%  "

# Row 40
$ws.Range("E40").Value = "  -1.22%  "

# Row 41
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.161"
$ws.Range("E41").Value = "  +4.16%  "

# Row 42
$ws.Range("B42").Value = "ImmutableX"
$ws.Range("C42").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.98"
$ws.Range("E42").Value = "  -1.13%  "

# Row 43
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.957"
$ws.Range("E43").Value = "  +2.68%  "

# Row 44
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "39.04"
$ws.Range("E44").Value = "  +17.92%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.47"
$ws.Range("E45").Value = "  +6.92%  "

# Row 46
$ws.Range("E46").Value = "  +0.04%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0465"
$ws.Range("E47").Value = "  +3.96%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.443"
$ws.Range("E48").Value = "  +10.25%  "

# Row 49
$ws.Range("E49").Value = "  -0.33%  "

# Row 50
$ws.Range("B50").Value = "MantraDAO"
$ws.Range("C50").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.69"
$ws.Range("E50").Value = "  +4.35%  "

# Row 51
$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.66"
$ws.Range("E51").Value = "  -0.09%  "

